$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.484.72'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '1.678.34'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '219.78'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').Value = '0.5345'
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.2713'
$ws.Range('E8').Value = '  +4.52%  '
$ws.Range('D9').Value = '0.06419'
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').Value = '21.97'
$ws.Range('E10').Value = '  +6.59%  '
$ws.Range('D11').Value = '0.07801'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.518'
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.677.54'
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').Value = '0.5602'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = '0.0₅8357'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '65.82'
$ws.Range('D17').Value = '26.529.17'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D19').Value = '4.813'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').Value = '193.76'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = '6.327'
$ws.Range('E22').Value = '  +2.81%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '142.34'
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').Value = '0.1286'
$ws.Range('E25').Value = '  +6.20%  '
$ws.Range('D26').Value = '7.425'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').Value = '16.33'
$ws.Range('E27').Value = '  +3.29%  '
$ws.Range('D28').Value = '1.443'
$ws.Range('E28').Value = '  +3.13%  '
$ws.Range('E29').Value = '  +5.89%  '
$ws.Range('E30').Value = '  +2.82%  '
$ws.Range('D31').Value = '3.610'
$ws.Range('E31').Value = '  +5.25%  '
$ws.Range('D32').Value = '3.465'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('D33').Value = '1.698'
$ws.Range('E33').Value = '  +3.58%  '
$ws.Range('D34').Value = '1.013'
$ws.Range('E34').Value = '  +3.21%  '
$ws.Range('D35').Value = '0.6165'
$ws.Range('E35').Value = '  +9.30%  '
$ws.Range('D36').Value = '2.421'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').Value = '2.787'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01636'
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '6.150'
$ws.Range('E39').Value = '  +7.87%  '
$ws.Range('D40').Value = '1.092.84'
$ws.Range('E40').Value = '  +5.93%  '
$ws.Range('D41').Value = '0.8679'
$ws.Range('E41').Value = '  +2.20%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '100.62'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = '1.824.20'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').Value = '0.0₈113'
$ws.Range('E45').Value = '  +2.47%  '
$ws.Range('D46').Value = '58.82'
$ws.Range('E46').Value = '  +5.46%  '
$ws.Range('D47').Value = '8.171'
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '6.053'
$ws.Range('E50').Value = '  +2.55%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.475'
$ws.Range('E51').Value = '  +6.85%  '
